$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLine = [char]10
$text15_16 = "Proyecto 2  + Manual Técnico " + $newLine + " ALVARO PEREZ NIÑO " + $newLine + " 703"
$text28 = "Proyecto 2  + Manual Técnico"

foreach ($col in @("B", "C", "D", "E")) {
    $ws.Range($col + "15").Value = $text15_16
    $ws.Range($col + "16").Value = $text15_16
}

$ws.Range("A28").Value = $text28
